# Sprint Planning - First week
# Populate sprint1 Sprint Backlog rows 4-7 with the first week's tasks,
# and switch the active sheet/selection from ListOfFeatures to sprint1.

$wb = $excel.ActiveWorkbook

$wsFeatures = $wb.Worksheets.Item("ListOfFeatures")
$wsSprint1  = $wb.Worksheets.Item("sprint1")

# --- Update sprint1 Sprint Backlog content (rows 4-7) ---

# Row 4
$wsSprint1.Range("A4").Value = "Task : Enroll to our 30 days free campaign! - db"
$wsSprint1.Range("B4").Value = 3
$wsSprint1.Range("C4").Value = "Mohan"
$wsSprint1.Range("D4").Value = "Ongoing"

# Row 5
$wsSprint1.Range("A5").Value = "Task : Contact Personal trainers online. -db"
$wsSprint1.Range("B5").Value = 3
$wsSprint1.Range("C5").Value = "Joni"
$wsSprint1.Range("D5").Value = "Ongoing"

# Row 6
$wsSprint1.Range("A6").Value = "Task : New payment methods  - db"
$wsSprint1.Range("B6").Value = 3
$wsSprint1.Range("C6").Value = "Edgar"
$wsSprint1.Range("D6").Value = "Ongoing"

# Row 7
$wsSprint1.Range("A7").Value = "Compare prices and pick one that suits you."
$wsSprint1.Range("B7").ClearContents()
$wsSprint1.Range("C7").Value = "All"
$wsSprint1.Range("D7").Value = "Ongoing"

# --- Switch active sheet/selection: sprint1 becomes active tab ---
# Set the (soon to be) inactive sheet's selection first, then activate
# sprint1 last so it ends up as the active/selected tab.
$wsFeatures.Range("A4").Select()

$wsSprint1.Activate()
$wsSprint1.Range("F23").Select()
